$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text that looks numeric (e.g. "597.13", "0.0301") as
# well as values using "." as a thousands separator (e.g. "67.397.08"). The
# source workbook stores these as plain text (inlineStr), so force the cell to
# text before writing and then drop back to the Normal style (no NumberFormat
# override survives) so no spurious style/format diff is introduced.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.397.08"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.88%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.493.08"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.03%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "597.13"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.60%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "180.39"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.60%  "

$ws.Range("B7").Value = "XRP"
$ws.Range("C7").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.607"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +4.01%  "

$ws.Range("B8").Value = "USDC"
$ws.Range("C8").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.489.20"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.06%  "

$ws.Range("E10").Value = "  +4.93%  "

$ws.Range("E11").Value = "  -1.69%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.437"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.60%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.093.68"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.07%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.33"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +9.74%  "

$ws.Range("E15").Value = "  +0.42%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "67.355.86"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.80%  "

$ws.Range("E17").Value = "  +0.08%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.493.38"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.50%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.30"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.61%  "

$ws.Range("E20").Value = "  +0.58%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "390.00"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.92%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.94"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.25%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "73.94"
$ws.Range("D23").Style = "Normal"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.543"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.78%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.998"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.34%  "

$ws.Range("E26").Value = "  +0.64%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000122"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.58%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.38"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.96%  "

$ws.Range("E29").Value = "  -2.79%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.57%  "

$ws.Range("E31").Value = "  +0.89%  "

$ws.Range("E32").Value = "  +0.28%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.07"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.02%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "23.54"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.23%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "7.39"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.71%  "

$ws.Range("E36").Value = "  +0.01%  "

$ws.Range("E37").Value = "  -0.30%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "163.22"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.36%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.871"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.62%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.81"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +11.29%  "

$ws.Range("E41").Value = "  -0.71%  "

$ws.Range("E42").Value = "  -0.42%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.65"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.45%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.851.95"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.42%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "26.44"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.69%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "26.83"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.16%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0723"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.87%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "41.67"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.28%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0301"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.36%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "333.24"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.49%  "

$ws.Range("E51").Value = "  -1.16%  "
